# Remove the "Id_membre" column (old column B); remaining columns shift left:
#   C->B (Prénom), D->C (Nom), E->D (Montant), F->E (Date enregistrement)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").EntireColumn.Delete()

# Drop the old 3rd data row (former row 4) - only two donation records remain
$ws.Range("A4").EntireRow.Delete()

# Header row: add the two new trailing columns
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Prénom"
$ws.Range("C1").Value = "Nom"
$ws.Range("D1").Value = "Montant"
$ws.Range("E1").Value = "Date enregistrement"
$ws.Range("F1").Value = "Statut"
$ws.Range("G1").Value = "Numéro de reçu"

# Row 2
$ws.Range("A2").Value = 204
$ws.Range("B2").Value = "celine"
$ws.Range("C2").Value = "celine"
$ws.Range("D2").Value = 9999
$ws.Range("E2").Value = "2019-06-21 12:31:58"
$ws.Range("F2").Value = "OK"
# Column G is left blank for this row (no receipt number yet)

# Row 3
$ws.Range("A3").Value = 157
$ws.Range("B3").Value = "Aazzouz"
$ws.Range("C3").Value = "Joelle"
$ws.Range("D3").Value = 5555
$ws.Range("E3").Value = "2019-06-05 09:20:36"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "R_2019_Foret_1001_Joelle_Aazzouz"
